# Attention-check rows inserted into the comprehension-questions table.
# Two new rows (new rows 6 & 7) are inserted before the former row 6,
# pushing the old rows 6-9 down to 8-11. The new rows ask participants
# to self-report how much attention they paid, with branching
# correct/incorrect feedback text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two blank rows at 6 & 7 (old rows 6-9 shift to 8-11) ---
$ws.Range("A6:A7").EntireRow.Insert()

# --- 2. Populate new row 6 (attention check on the images) ---
$ws.Range("A6").Value = "Mit gondol, ebben a körben milyen arányban tudott figyelni a képekre? Kérjük, válaszoljon őszintén. Válasza a vizsgálat bejefezésével járó jutalom (kredit, ajándékutalvány) értékét nem befolyásolja."
$ws.Range("B6").Value = "D: 0-24%    F: 25-49%    J: 50-74%    K: 75-100%"
$ws.Range("C6").Value = "j"
$ws.Range("D6").Value = "Köszönjük!"
$ws.Range("E6").Value = "Ön ennek a körnek kevesebb, mint háromnegyedében tudott a képekre figyelni."
$ws.Range("F6").Value = "Kérjük, próbáljon meg a továbbiakban figyelni a képekre. "

# --- 3. Populate new row 7 (attention check on the image locations) ---
$ws.Range("A7").Value = "Mit gondol, ebben a körben milyen arányban tudott figyelni a képek helyére? Kérjük, válaszoljon őszintén. Válasza a vizsgálat bejefezésével járó jutalom (kredit, ajándékutalvány) értékét nem befolyásolja."
$ws.Range("B7").Value = "D: 0-24%    F: 25-49%    J: 50-74%    K: 75-100%"
$ws.Range("C7").Value = "j"
$ws.Range("D7").Value = "Köszönjük!"
$ws.Range("E7").Value = "Ön ebben a körnek kevesebb, mint háromnegyedében tudott a képek helyszínére figyelni."
$ws.Range("F7").Value = "Kérjük, próbáljon meg a továbbiakban figyelni a képek helyszínére. "

# --- 4. Row heights ---
$ws.Rows(4).RowHeight = 181.5
$ws.Rows(5).RowHeight = 181.5
$ws.Rows(6).RowHeight = 35.05
$ws.Rows(7).RowHeight = 35.05
$ws.Rows(8).RowHeight = 125.25
$ws.Rows(9).RowHeight = 125.25

# --- 5. Wrap text on the new rows' long-text cells ---
$ws.Range("A6:A7").WrapText = $true
$ws.Range("E6:F7").WrapText = $true

# --- 6. Column E gets its own (non-wrapping) style, split off from column F ---
$ws.Range("E1:E9").WrapText = $false
$ws.Columns("E").ColumnWidth = $ws.Columns("F").ColumnWidth

# --- 7. Sheet view: zoom out a bit and focus the view near the new rows ---
$ws.Application.ActiveWindow.Zoom = 90
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("E8").Select()
